# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型
# sheets to reflect the latest scrape (gh-pages output regenerated at
# commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 200
$wsExpo.Range("F6").Value = 507
$wsExpo.Range("F7").Value = 106
$wsExpo.Range("F8").Value = 120
$wsExpo.Range("F9").Value = 48
$wsExpo.Range("F10").Value = 6794
$wsExpo.Range("F12").Value = 378
$wsExpo.Range("F13").Value = 3109
$wsExpo.Range("F14").Value = 205
$wsExpo.Range("F15").Value = 359
$wsExpo.Range("F16").Value = 260
$wsExpo.Range("F17").Value = 550
$wsExpo.Range("F18").Value = 8

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 200
$wsAll.Range("F8").Value = 507
$wsAll.Range("F9").Value = 106
$wsAll.Range("F10").Value = 120
$wsAll.Range("F11").Value = 48
$wsAll.Range("F13").Value = 6794
$wsAll.Range("F16").Value = 378
$wsAll.Range("F17").Value = 3109
$wsAll.Range("F18").Value = 205
$wsAll.Range("F19").Value = 359
$wsAll.Range("F20").Value = 260
$wsAll.Range("F21").Value = 550
$wsAll.Range("F22").Value = 8
